$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.497.89'
$ws.Range("E2").Value = '  -0.76%  '
$ws.Range("D3").Value = '1.825.60'
$ws.Range("E3").Value = '  -1.36%  '
$ws.Range("E4").Value = '  -0.33%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '312.39'
$ws.Range("E5").Value = '  -0.15%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  -0.43%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4253'
$ws.Range("E7").Value = '  -0.64%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3612'
$ws.Range("E8").Value = '  +0.76%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07208'
$ws.Range("E9").Value = '  -1.28%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8619'
$ws.Range("E10").Value = '  -1.10%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.63'
$ws.Range("E11").Value = '  -0.54%  '
$ws.Range("D12").Value = '1.799.18'
$ws.Range("E12").Value = '  -2.70%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.393'
$ws.Range("E13").Value = '  +1.07%  '
$ws.Range("E14").Value = '  -1.11%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.06923'
$ws.Range("E15").Value = '  -1.45%  '
$ws.Range("E16").Value = '  -0.58%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '80.60'
$ws.Range("E17").Value = '  +1.20%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008923'
$ws.Range("E18").Value = '  -0.34%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.002'
$ws.Range("E19").Value = '  -0.23%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.36'
$ws.Range("E20").Value = '  +0.41%  '
$ws.Range("D21").Value = '27.550.09'
$ws.Range("E21").Value = '  -0.55%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.119'
$ws.Range("E22").Value = '  +2.31%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.89'
$ws.Range("E23").Value = '  +4.81%  '
$ws.Range("D24").Value = '2.058.65'
$ws.Range("E24").Value = '  -0.78%  '
$ws.Range("E25").Value = '  +0.06%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '154.98'
$ws.Range("E26").Value = '  -0.61%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.71'
$ws.Range("E27").Value = '  +1.02%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.168'
$ws.Range("E28").Value = '  -2.06%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '114.13'
$ws.Range("E29").Value = '  -5.43%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.794'
$ws.Range("E30").Value = '  -4.22%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08875'
$ws.Range("E31").Value = '  -0.51%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7512'
$ws.Range("E32").Value = '  -0.90%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.543'
$ws.Range("E33").Value = '  +0.65%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.972'
$ws.Range("E34").Value = '  -0.08%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.124'
$ws.Range("E35").Value = '  +0.04%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.000'
$ws.Range("E36").Value = '  -0.39%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.088'
$ws.Range("E37").Value = '  -1.21%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05291'
$ws.Range("E38").Value = '  -2.59%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01922'
$ws.Range("E39").Value = '  -0.47%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.793'
$ws.Range("E40").Value = '  -1.22%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5075'
$ws.Range("E41").Value = '  -0.04%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1657'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.459'
$ws.Range("E43").Value = '  -2.22%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.358'
$ws.Range("E44").Value = '  -0.55%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.49'
$ws.Range("E45").Value = '  +1.12%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '106.20'
$ws.Range("E46").Value = '  -0.03%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.06467'
$ws.Range("E47").Value = '  -1.12%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4678'
$ws.Range("E48").Value = '  +0.28%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.9995'
$ws.Range("E49").Value = '  -0.55%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.617'
$ws.Range("E50").Value = '  -0.77%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '63.84'
$ws.Range("E51").Value = '  -1.04%  '
